$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three sensor rows that are being grouped/consolidated away.
# Delete from the bottom up so earlier row numbers remain valid targets.
$ws.Rows.Item(66).Delete()
$ws.Rows.Item(36).Delete()
$ws.Rows.Item(21).Delete()

# Update the selected cell to match the saved view state.
$ws.Range("I12").Select()
